$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 96 - duplicate of row 95's layout with new employee data
$ws.Range("C95:F95").Copy()
$ws.Range("C96").PasteSpecial(-4122)

$ws.Range("A96").Value = 201
$ws.Range("B96").Value = "Test Karna"
$ws.Range("C96").Value = $ws.Range("C95").Value()
$ws.Range("D96").Value = $ws.Range("D95").Value()
$ws.Range("E96").Value = $ws.Range("E95").Value()
$ws.Range("F96").Value = $ws.Range("F95").Value()

# Clear the footer text cells (A97, C97) - the "Generated on" / "report generated by" strings are removed
$ws.Range("A97").Value = ""
$ws.Range("C97").Value = ""

# Update the active selection to reflect where the user ended up after the edit
[void]$ws.Range("D101").Select()
